$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the date formatting from the cell above (B6) so the new Date
# cells (B7, B8) reuse the existing "short date" style instead of
# Excel creating a brand-new number format / style entry.
$ws.Range("B6").Copy()
$ws.Range("B7:B8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Week 7 (row 7) ---
$ws.Range("B7").Value = 45355          # Date for week 7 (2024-03-04)
$ws.Range("C7").Value = 0.166
$ws.Range("D7").Value = 0.166
$ws.Range("E7").Value = 0.166
$ws.Range("F7").Value = 0.166
$ws.Range("G7").Value = 0.166
$ws.Range("H7").Value = 0.166

# --- Week 8 (row 8) ---
$ws.Range("B8").Value = 45362          # Date for week 8 (2024-03-11)
$ws.Range("C8").Value = 0.166
$ws.Range("D8").Value = 0.166
$ws.Range("E8").Value = 0.166
$ws.Range("F8").Value = 0.166
$ws.Range("G8").Value = 0.166
$ws.Range("H8").Value = 0.166

# Make sure dependent formulas (J7, J8, J15, etc.) are recalculated
$excel.Calculate()

# --- Update the sheet selection to match the saved state ---
$ws.Range("H8").Select()
